# "Status by State" pivot sheet: the "New" column (C) and the "Present"
# column (D) are being consolidated into a single column. For every data
# row, the new value of column C becomes the old "New" value plus the old
# "Present" value, the header of column C becomes "Present" (taking over
# the old header of column D), and column D is then removed entirely so
# the sheet's used range shrinks from A1:D39 to A1:C39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status by State")

$lastRow = 39

for ($r = 2; $r -le $lastRow; $r++) {
    $newVal = $ws.Cells.Item($r, 3).Value()
    $presentVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $newVal + $presentVal
}

$ws.Cells.Item(1, 3).Value = "Present"

$ws.Columns.Item(4).EntireColumn.Delete()
